# The source submission file had "sex" recorded as the abbreviation "F";
# this updates it to the full word "female" (value lives in B2, under the
# "sex" header in B1) and moves the active selection to B2 to match the
# session state captured when the edit was made.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "female"
$ws.Range("B2").Select()
